# fix: fix bug of WireBuilder dead lock
#
# The WireBuilder tool (which generates these "uart"/"uart_rx"/"uart_tx"
# wiring-table sheets) had ordered the port rows on the top-level "uart"
# sheet in a way that caused a dead lock; the fix re-orders the rows on
# that sheet (grouping the just-added rx-side wires together) and fills
# in a few Port-info annotations that the tool had left behind / stale.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "uart": reorder the port rows 8-15 and annotate Port-info (col E)
# ---------------------------------------------------------------------
$wsUart = $wb.Worksheets.Item("uart")

# Re-write rows 8..15 (Port-name / InOut) in the new order. Width (col C)
# is always 1 for all of these rows, so it is left untouched.
$wsUart.Range("A8").Value  = "rxd"
$wsUart.Range("B8").Value  = "input"

$wsUart.Range("A9").Value  = "s_axis_tready"
$wsUart.Range("B9").Value  = "output"

$wsUart.Range("A10").Value = "s_axis_tvalid"
$wsUart.Range("B10").Value = "input"

$wsUart.Range("A11").Value = "txd"
$wsUart.Range("B11").Value = "output"

$wsUart.Range("A12").Value = "rx_busy"
$wsUart.Range("B12").Value = "output"

$wsUart.Range("A13").Value = "rx_frame_error"
$wsUart.Range("B13").Value = "output"

$wsUart.Range("A14").Value = "rx_overrun_error"
$wsUart.Range("B14").Value = "output"

$wsUart.Range("A15").Value = "tx_busy"
$wsUart.Range("B15").Value = "output"

# New Port-info notes
$wsUart.Range("E5").Value = "fsdf"
$wsUart.Range("E7").Value = "sfdf"
$wsUart.Range("E9").Value = "fsdf"

# ---------------------------------------------------------------------
# Sheet "uart_rx": drop the stale Port-info note on "rxd" and move/refresh
# the notes attached to "m_axis_tdata" / "overrun_error"
# ---------------------------------------------------------------------
$wsRx = $wb.Worksheets.Item("uart_rx")

$wsRx.Range("E5").Value  = "sdf"
$wsRx.Range("E8").Value  = ""
$wsRx.Range("E10").Value = "rx_overrun_error"

# ---------------------------------------------------------------------
# Sheet "uart_tx": drop the stale Port-info note on "s_axis_tdata" and
# add fresh notes to "rst" / "s_axis_tready"
# ---------------------------------------------------------------------
$wsTx = $wb.Worksheets.Item("uart_tx")

$wsTx.Range("E4").Value = "ert"
$wsTx.Range("E5").Value = ""
$wsTx.Range("E7").Value = "dfg"
